{"js": "// Replace the two-digit multiplication answers in the body text.\n// Each old equation string is unique in the document, so an exact,\n// case-sensitive search-and-replace is unambiguous for every cell.\nconst replacements = [\n  [\"69\u00d794=6486\", \"94\u00d788=8272\"],\n  [\"85\u00d723=1955\", \"33\u00d795=3135\"],\n  [\"87\u00d771=6177\", \"38\u00d747=1786\"],\n  [\"19\u00d765=1235\", \"16\u00d721=336\"],\n  [\"29\u00d758=1682\", \"69\u00d737=2553\"],\n  [\"89\u00d738=3382\", \"66\u00d743=2838\"],\n  [\"24\u00d778=1872\", \"35\u00d787=3045\"],\n  [\"53\u00d734=1802\", \"55\u00d728=1540\"],\n  [\"12\u00d768=816\", \"52\u00d757=2964\"],\n  [\"75\u00d778=5850\", \"86\u00d730=2580\"],\n  [\"94\u00d767=6298\", \"65\u00d721=1365\"],\n  [\"81\u00d726=2106\", \"11\u00d737=407\"],\n  [\"37\u00d751=1887\", \"49\u00d783=4067\"],\n  [\"77\u00d790=6930\", \"49\u00d791=4459\"],\n  [\"36\u00d756=2016\", \"50\u00d789=4450\"],\n  [\"87\u00d749=4263\", \"29\u00d771=2059\"],\n  [\"35\u00d712=420\", \"96\u00d798=9408\"],\n  [\"81\u00d757=4617\", \"94\u00d712=1128\"],\n  [\"55\u00d787=4785\", \"40\u00d741=1640\"],\n  [\"31\u00d779=2449\", \"68\u00d738=2584\"],\n  [\"46\u00d796=4416\", \"26\u00d759=1534\"],\n  [\"34\u00d778=2652\", \"72\u00d765=4680\"],\n  [\"92\u00d725=2300\", \"30\u00d730=900\"],\n  [\"41\u00d745=1845\", \"64\u00d712=768\"],\n  [\"95\u00d778=7410\", \"35\u00d743=1505\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit multiplication answers in the body text.\n# Each old equation string is unique in the document, so an exact,\n# case-sensitive Find/Replace is unambiguous for every table cell.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"69\u00d794=6486\", \"94\u00d788=8272\"),\n    @(\"85\u00d723=1955\", \"33\u00d795=3135\"),\n    @(\"87\u00d771=6177\", \"38\u00d747=1786\"),\n    @(\"19\u00d765=1235\", \"16\u00d721=336\"),\n    @(\"29\u00d758=1682\", \"69\u00d737=2553\"),\n    @(\"89\u00d738=3382\", \"66\u00d743=2838\"),\n    @(\"24\u00d778=1872\", \"35\u00d787=3045\"),\n    @(\"53\u00d734=1802\", \"55\u00d728=1540\"),\n    @(\"12\u00d768=816\", \"52\u00d757=2964\"),\n    @(\"75\u00d778=5850\", \"86\u00d730=2580\"),\n    @(\"94\u00d767=6298\", \"65\u00d721=1365\"),\n    @(\"81\u00d726=2106\", \"11\u00d737=407\"),\n    @(\"37\u00d751=1887\", \"49\u00d783=4067\"),\n    @(\"77\u00d790=6930\", \"49\u00d791=4459\"),\n    @(\"36\u00d756=2016\", \"50\u00d789=4450\"),\n    @(\"87\u00d749=4263\", \"29\u00d771=2059\"),\n    @(\"35\u00d712=420\", \"96\u00d798=9408\"),\n    @(\"81\u00d757=4617\", \"94\u00d712=1128\"),\n    @(\"55\u00d787=4785\", \"40\u00d741=1640\"),\n    @(\"31\u00d779=2449\", \"68\u00d738=2584\"),\n    @(\"46\u00d796=4416\", \"26\u00d759=1534\"),\n    @(\"34\u00d778=2652\", \"72\u00d765=4680\"),\n    @(\"92\u00d725=2300\", \"30\u00d730=900\"),\n    @(\"41\u00d745=1845\", \"64\u00d712=768\"),\n    @(\"95\u00d778=7410\", \"35\u00d743=1505\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $ok = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $ok) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
